$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (ID) values must be stored as text that looks like numbers.
# Temporarily apply a text number format so Excel doesn't auto-convert the
# values to numeric cells, then clear the formatting again so the cells end
# up with no explicit style (matching a plain copy/paste of text values).
$idRange = $ws.Range("A2:A8")
$idRange.NumberFormat = "@"

$ws.Range("A2").Value = "20"
$ws.Range("A3").Value = "10"
$ws.Range("A4").Value = "23"
$ws.Range("A5").Value = "33"
$ws.Range("A6").Value = "5"
$ws.Range("A7").Value = "0"
$ws.Range("A8").Value = "23"

$idRange.ClearFormats()

# Column B (Value) values are plain numbers.
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 30
$ws.Range("B5").Value = 53
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 100
